$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"


# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("M15").Value = 300
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 108
$ws.Range("J16").Value = 128
$ws.Range("K16").Value = -15.625
$ws.Range("L16").Value = 61.194029850746
$ws.Range("M16").Value = 4.854368932038
$ws.Range("N16").Value = -81.909547738693
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 35.294117647058
$ws.Range("I17").Value = 159
$ws.Range("J17").Value = 173
$ws.Range("K17").Value = -8.092485549132
$ws.Range("L17").Value = 8.163265306122
$ws.Range("M17").Value = 133.823529411765
$ws.Range("N17").Value = -28.054298642533
$ws.Range("C18").Value = 12
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 148
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = -7.5
$ws.Range("L18").Value = 29.824561403508
$ws.Range("M18").Value = -12.941176470588
$ws.Range("N18").Value = -89.664804469273
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -8.620689655172
$ws.Range("I19").Value = 456
$ws.Range("J19").Value = 518
$ws.Range("K19").Value = -11.969111969112
$ws.Range("L19").Value = 68.888888888888
$ws.Range("M19").Value = 80.952380952380
$ws.Range("N19").Value = -53.846153846153
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 27.777777777777
$ws.Range("I20").Value = 140
$ws.Range("J20").Value = 128
$ws.Range("K20").Value = 9.375
$ws.Range("L20").Value = 169.230769230769
$ws.Range("M20").Value = 1.449275362318
$ws.Range("N20").Value = -93.889131383675
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = 15.447154471544
$ws.Range("I21").Value = 1027
$ws.Range("J21").Value = 1114
$ws.Range("K21").Value = -7.809694793536
$ws.Range("L21").Value = 56.793893129771
$ws.Range("M21").Value = 39.348710990502
$ws.Range("N21").Value = -81.485487650982
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 233.333333333333
$ws.Range("M22").Value = 11.111111111111
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 11
$ws.Range("I23").Value = 45
$ws.Range("J23").Value = 65
$ws.Range("K23").Value = -30.769230769230
$ws.Range("L23").Value = -2.173913043478
$ws.Range("M23").Value = 164.705882352941
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -38.461538461538
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 161
$ws.Range("H24").Value = -38.509316770186
$ws.Range("I24").Value = 831
$ws.Range("J24").Value = 1084
$ws.Range("K24").Value = -23.339483394833
$ws.Range("L24").Value = 23.476968796433
$ws.Range("M24").Value = 16.223776223776
$ws.Range("C25").Value = 8
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 6.896551724137
$ws.Range("I25").Value = 273
$ws.Range("J25").Value = 305
$ws.Range("K25").Value = -10.491803278688
$ws.Range("L25").Value = 20.264317180616
$ws.Range("M25").Value = -19.705882352941
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = 7.142857142857
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = -6.451612903225
$ws.Range("H30").Value = -100

# --- Cells changing between numeric and text (N/A / ***.* placeholders) ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 4

$ws.Range("K14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = -50

$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
